$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 275; this shifts existing rows 275-335 down to 276-336
# (and the sheet dimension grows from A1:R335 to A1:R336), reproducing the
# weekly roll described in the diff where every prior row's daily record
# moves down one slot and a brand-new "latest" record is prepended.
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new record.
$ws.Range("A275").Value = 10
$ws.Range("B275").Value = "Vega Modelo de Temuco"
$ws.Range("C275").Value = "La Araucanía"
$ws.Range("D275").Value = 44511
$ws.Range("E275").Value = 9
$ws.Range("F275").Value = 100112028
$ws.Range("G275").Value = "Sandia"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 800
$ws.Range("K275").Value = 700
$ws.Range("L275").Value = 700
$ws.Range("M275").Value = 700
$ws.Range("N275").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O275").Value = "Perú"
$ws.Range("P275").Value = 700
$ws.Range("Q275").Value = 1
$ws.Range("R275").Value = "Hortaliza"
